# Auto-generated script: apply 2023-03-06 daily crime count increments
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1118
$ws.Range("I3").Value = 7488
$ws.Range("J3").Value = 1199
$ws.Range("J5").Value = 1312
$ws.Range("F6").Value = 1880
$ws.Range("I6").Value = 1752
$ws.Range("J6").Value = 263
$ws.Range("J7").Value = 90
$ws.Range("J8").Value = 5342
$ws.Range("I9").Value = 8968
$ws.Range("J9").Value = 1596
$ws.Range("I10").Value = 54672
$ws.Range("J10").Value = 8573
$ws.Range("F11").Value = 105547
$ws.Range("I11").Value = 110300
$ws.Range("J11").Value = 19574

$ws = $wb.Sheets.Item('By Neighborhood')
$ws.Range("J2").Value = 183
$ws.Range("J4").Value = 98
$ws.Range("J6").Value = 196
$ws.Range("J7").Value = 506
$ws.Range("J8").Value = 808
$ws.Range("J9").Value = 104
$ws.Range("J10").Value = 172
$ws.Range("J11").Value = 322
$ws.Range("J15").Value = 172
$ws.Range("J16").Value = 130
$ws.Range("J17").Value = 27
$ws.Range("J18").Value = 149
$ws.Range("J19").Value = 510
$ws.Range("J20").Value = 326
$ws.Range("J22").Value = 88
$ws.Range("J23").Value = 269
$ws.Range("J24").Value = 109
$ws.Range("J25").Value = 85
$ws.Range("J27").Value = 238
$ws.Range("J29").Value = 708
$ws.Range("J33").Value = 488
$ws.Range("J34").Value = 177
$ws.Range("J36").Value = 260
$ws.Range("J37").Value = 518
$ws.Range("J41").Value = 95
$ws.Range("J42").Value = 523
$ws.Range("J43").Value = 212
$ws.Range("J46").Value = 70
$ws.Range("J47").Value = 171
$ws.Range("J48").Value = 407
$ws.Range("J49").Value = 240
$ws.Range("J50").Value = 192
$ws.Range("J51").Value = 264
$ws.Range("J52").Value = 291
$ws.Range("J53").Value = 308
$ws.Range("J54").Value = 582
$ws.Range("J56").Value = 94
$ws.Range("J60").Value = 140
$ws.Range("F63").Value = 1316
$ws.Range("I63").Value = 2253
$ws.Range("J63").Value = 220
$ws.Range("J64").Value = 168
$ws.Range("J65").Value = 291
$ws.Range("I67").Value = 2510
$ws.Range("J67").Value = 421
$ws.Range("J68").Value = 63
$ws.Range("J69").Value = 75
$ws.Range("J70").Value = 131
$ws.Range("J72").Value = 97
$ws.Range("J73").Value = 213
$ws.Range("J75").Value = 75
$ws.Range("J76").Value = 575
$ws.Range("J77").Value = 94
$ws.Range("J78").Value = 294
$ws.Range("J79").Value = 440
$ws.Range("J80").Value = 56
$ws.Range("J83").Value = 350
$ws.Range("J84").Value = 163
$ws.Range("J85").Value = 762
$ws.Range("J86").Value = 132
$ws.Range("J87").Value = 73
$ws.Range("J88").Value = 169
$ws.Range("J89").Value = 288
$ws.Range("J91").Value = 196
$ws.Range("J93").Value = 145
$ws.Range("J94").Value = 413
$ws.Range("J95").Value = 284
$ws.Range("J96").Value = 281
$ws.Range("J97").Value = 290
$ws.Range("J98").Value = 177
$ws.Range("J99").Value = 276
$ws.Range("F101").Value = 105547
$ws.Range("I101").Value = 110300
$ws.Range("J101").Value = 19574

$ws = $wb.Sheets.Item('Uptown')
$ws.Range("J6").Value = 4
$ws.Range("J8").Value = 61
$ws.Range("J10").Value = 152
$ws.Range("J11").Value = 288

$ws = $wb.Sheets.Item('South Shore')
$ws.Range("J2").Value = 44
$ws.Range("J3").Value = 63
$ws.Range("J8").Value = 288
$ws.Range("J9").Value = 49
$ws.Range("J10").Value = 224
$ws.Range("J11").Value = 762

$ws = $wb.Sheets.Item('Norwood Park')
$ws.Range("J3").Value = 3
$ws.Range("J11").Value = 75

$ws = $wb.Sheets.Item('Little Village')
$ws.Range("J9").Value = 37
$ws.Range("J10").Value = 127
$ws.Range("J11").Value = 291

$ws = $wb.Sheets.Item('Belmont Cragin')
$ws.Range("J7").Value = 3
$ws.Range("J8").Value = 94
$ws.Range("J10").Value = 151
$ws.Range("J11").Value = 322

$ws = $wb.Sheets.Item('Austin')
$ws.Range("J2").Value = 89
$ws.Range("J3").Value = 88
$ws.Range("J8").Value = 251
$ws.Range("J10").Value = 246
$ws.Range("J11").Value = 808

$ws = $wb.Sheets.Item('Logan Square')
$ws.Range("J2").Value = 5
$ws.Range("J8").Value = 99
$ws.Range("J9").Value = 24
$ws.Range("J10").Value = 143
$ws.Range("J11").Value = 308

$ws = $wb.Sheets.Item('Auburn Gresham')
$ws.Range("J5").Value = 52
$ws.Range("J8").Value = 173
$ws.Range("J10").Value = 156
$ws.Range("J11").Value = 506

$ws = $wb.Sheets.Item('West Ridge')
$ws.Range("J2").Value = 19
$ws.Range("J5").Value = 29
$ws.Range("J6").Value = 4
$ws.Range("J10").Value = 128
$ws.Range("J11").Value = 281

$ws = $wb.Sheets.Item('O''Hare')
$ws.Range("J10").Value = 99
$ws.Range("J11").Value = 131

$ws = $wb.Sheets.Item('Grand Crossing')
$ws.Range("J2").Value = 39
$ws.Range("J5").Value = 32
$ws.Range("J8").Value = 192
$ws.Range("J10").Value = 145
$ws.Range("J11").Value = 518

$ws = $wb.Sheets.Item('Woodlawn')
$ws.Range("J3").Value = 18
$ws.Range("J8").Value = 102
$ws.Range("J10").Value = 87
$ws.Range("J11").Value = 276

$ws = $wb.Sheets.Item('North Lawndale')
$ws.Range("I3").Value = 365
$ws.Range("J3").Value = 62
$ws.Range("J8").Value = 112
$ws.Range("J9").Value = 52
$ws.Range("J10").Value = 125
$ws.Range("I11").Value = 2510
$ws.Range("J11").Value = 421

$ws = $wb.Sheets.Item('Gage Park')
$ws.Range("I5").Value = 67
$ws.Range("I9").Value = 101

$ws = $wb.Sheets.Item('South Deering')
$ws.Range("J8").Value = 71
$ws.Range("J11").Value = 163

$ws = $wb.Sheets.Item('New City')
$ws.Range("J2").Value = 27
$ws.Range("J10").Value = 98
$ws.Range("J11").Value = 291

$ws = $wb.Sheets.Item('Edgewater')
$ws.Range("J5").Value = 33
$ws.Range("J10").Value = 143
$ws.Range("J11").Value = 238

$ws = $wb.Sheets.Item('South Chicago')
$ws.Range("J8").Value = 127
$ws.Range("J11").Value = 350

$ws = $wb.Sheets.Item('Garfield Park')
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 48
$ws.Range("J8").Value = 140
$ws.Range("J10").Value = 157
$ws.Range("J11").Value = 488

$ws = $wb.Sheets.Item('West Pullman')
$ws.Range("J5").Value = 13
$ws.Range("J8").Value = 111
$ws.Range("J9").Value = 23
$ws.Range("J11").Value = 284

$ws = $wb.Sheets.Item('Pullman')
$ws.Range("J8").Value = 25
$ws.Range("J11").Value = 75

$ws = $wb.Sheets.Item('Roseland')
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 42
$ws.Range("J8").Value = 168
$ws.Range("J10").Value = 119
$ws.Range("J11").Value = 440

$ws = $wb.Sheets.Item('Jefferson Park')
$ws.Range("J10").Value = 29
$ws.Range("J11").Value = 70

$ws = $wb.Sheets.Item('Morgan Park')
$ws.Range("J3").Value = 8
$ws.Range("J10").Value = 55
$ws.Range("J11").Value = 140

$ws = $wb.Sheets.Item('Near South Side')
$ws.Range("J8").Value = 54
$ws.Range("J10").Value = 79
$ws.Range("J11").Value = 168

$ws = $wb.Sheets.Item('West Loop')
$ws.Range("J8").Value = 32
$ws.Range("J10").Value = 337
$ws.Range("J11").Value = 413

$ws = $wb.Sheets.Item('River North')
$ws.Range("J2").Value = 10
$ws.Range("J3").Value = 17
$ws.Range("J8").Value = 49
$ws.Range("J10").Value = 429
$ws.Range("J11").Value = 575

$ws = $wb.Sheets.Item('Ukrainian Village')
$ws.Range("J9").Value = 10
$ws.Range("J11").Value = 73

$ws = $wb.Sheets.Item('East Side')
$ws.Range("J8").Value = 13
$ws.Range("J9").Value = 5
$ws.Range("J10").Value = 41
$ws.Range("J11").Value = 85

$ws = $wb.Sheets.Item('Bucktown')
$ws.Range("J10").Value = 85
$ws.Range("J11").Value = 130

$ws = $wb.Sheets.Item('Lincoln Park')
$ws.Range("J10").Value = 167
$ws.Range("J11").Value = 240

$ws = $wb.Sheets.Item('West Town')
$ws.Range("J5").Value = 19
$ws.Range("J11").Value = 290

$ws = $wb.Sheets.Item('Loop')
$ws.Range("J3").Value = 13
$ws.Range("J8").Value = 70
$ws.Range("J9").Value = 41
$ws.Range("J10").Value = 422
$ws.Range("J11").Value = 582

$ws = $wb.Sheets.Item('Portage Park')
$ws.Range("J2").Value = 15
$ws.Range("J8").Value = 46
$ws.Range("J10").Value = 109
$ws.Range("J11").Value = 213

$ws = $wb.Sheets.Item('Englewood')
$ws.Range("J2").Value = 67
$ws.Range("J8").Value = 198
$ws.Range("J9").Value = 63
$ws.Range("J11").Value = 708

$ws = $wb.Sheets.Item('Chatham')
$ws.Range("J8").Value = 159
$ws.Range("J9").Value = 58
$ws.Range("J10").Value = 166
$ws.Range("J11").Value = 510

$ws = $wb.Sheets.Item('Clearing')
$ws.Range("J2").Value = 3
$ws.Range("J11").Value = 88

$ws = $wb.Sheets.Item('Humboldt Park')
$ws.Range("J8").Value = 151
$ws.Range("J9").Value = 107
$ws.Range("J10").Value = 166
$ws.Range("J11").Value = 523

$ws = $wb.Sheets.Item('Lake View')
$ws.Range("I5").Value = 180
$ws.Range("I10").Value = 1773
$ws.Range("J10").Value = 280
$ws.Range("J11").Value = 407

$ws = $wb.Sheets.Item('Ashburn')
$ws.Range("J8").Value = 71
$ws.Range("J10").Value = 70
$ws.Range("J11").Value = 196

$ws = $wb.Sheets.Item('Hermosa')
$ws.Range("J3").Value = 6
$ws.Range("J8").Value = 41
$ws.Range("J10").Value = 24
$ws.Range("J11").Value = 95

$ws = $wb.Sheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 19
$ws.Range("J10").Value = 85
$ws.Range("J11").Value = 260

$ws = $wb.Sheets.Item('Avondale')
$ws.Range("J8").Value = 37
$ws.Range("J10").Value = 93
$ws.Range("J11").Value = 172

$ws = $wb.Sheets.Item('Streeterville')
$ws.Range("J10").Value = 90
$ws.Range("J11").Value = 132

$ws = $wb.Sheets.Item('Rogers Park')
$ws.Range("J10").Value = 156
$ws.Range("J11").Value = 294

$ws = $wb.Sheets.Item('North Park')
$ws.Range("J8").Value = 11
$ws.Range("J11").Value = 63

$ws = $wb.Sheets.Item('Brighton Park')
$ws.Range("J2").Value = 15
$ws.Range("J5").Value = 13
$ws.Range("J10").Value = 65
$ws.Range("J11").Value = 172

$ws = $wb.Sheets.Item('Dunning')
$ws.Range("J5").Value = 17
$ws.Range("J11").Value = 109

$ws = $wb.Sheets.Item('Douglas')
$ws.Range("J10").Value = 101
$ws.Range("J11").Value = 269

$ws = $wb.Sheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 33
$ws.Range("J10").Value = 139
$ws.Range("J11").Value = 326

$ws = $wb.Sheets.Item('Little Italy, UIC')
$ws.Range("J8").Value = 95
$ws.Range("J10").Value = 99
$ws.Range("J11").Value = 264

$ws = $wb.Sheets.Item('Washington Park')
$ws.Range("J3").Value = 25
$ws.Range("J8").Value = 80
$ws.Range("J10").Value = 49
$ws.Range("J11").Value = 196

$ws = $wb.Sheets.Item('Kenwood')
$ws.Range("J10").Value = 70
$ws.Range("J11").Value = 171

$ws = $wb.Sheets.Item('Lincoln Square')
$ws.Range("J8").Value = 43
$ws.Range("J10").Value = 86
$ws.Range("J11").Value = 192

$ws = $wb.Sheets.Item('West Lawn')
$ws.Range("J8").Value = 51
$ws.Range("J11").Value = 145

$ws = $wb.Sheets.Item('Calumet Heights')
$ws.Range("J2").Value = 15
$ws.Range("J8").Value = 42
$ws.Range("J11").Value = 149

$ws = $wb.Sheets.Item('Riverdale')
$ws.Range("J7").Value = 2
$ws.Range("J8").Value = 33
$ws.Range("J11").Value = 94

$ws = $wb.Sheets.Item('Magnificent Mile')
$ws.Range("J10").Value = 87
$ws.Range("J11").Value = 94

$ws = $wb.Sheets.Item('Albany Park')
$ws.Range("J10").Value = 98
$ws.Range("J11").Value = 183

$ws = $wb.Sheets.Item('Old Town')
$ws.Range("J10").Value = 52
$ws.Range("J11").Value = 97

$ws = $wb.Sheets.Item('Hyde Park')
$ws.Range("J6").Value = 3
$ws.Range("J10").Value = 101
$ws.Range("J11").Value = 212

$ws = $wb.Sheets.Item('Burnside')
$ws.Range("J2").Value = 2
$ws.Range("J11").Value = 27

$ws = $wb.Sheets.Item('Archer Heights')
$ws.Range("J9").Value = 8
$ws.Range("J10").Value = 41
$ws.Range("J11").Value = 98

$ws = $wb.Sheets.Item('Rush & Division')
$ws.Range("J10").Value = 39
$ws.Range("J11").Value = 56

$ws = $wb.Sheets.Item('Garfield Ridge')
$ws.Range("J8").Value = 53
$ws.Range("J10").Value = 79
$ws.Range("J11").Value = 177

$ws = $wb.Sheets.Item('Wicker Park')
$ws.Range("J10").Value = 94
$ws.Range("J11").Value = 177

$ws = $wb.Sheets.Item('Avalon Park')
$ws.Range("J8").Value = 44
$ws.Range("J11").Value = 104

$ws = $wb.Sheets.Item('United Center')
$ws.Range("J8").Value = 62
$ws.Range("J10").Value = 65
$ws.Range("J11").Value = 169
